$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$iValues = @(6,9,5,8,6,7,6,7,6,5,7,8,9,6,5,7,4,8,7,8)
$jValues = @(6,9,5,8,7,7,6,7,6,5,7,9,9,6,6,7,5,8,7,8)

for ($r = 0; $r -lt $iValues.Length; $r++) {
    $row = $r + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$r]
    $ws.Cells.Item($row, 10).Value = $jValues[$r]
}
